$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-26 Tuesday", "2024-03-27 Wednesday"),
    @("984÷4=", "118÷5="),
    @("298÷4=", "769÷2="),
    @("368÷7=", "941÷5="),
    @("936÷4=", "999÷8="),
    @("830÷7=", "251÷4="),
    @("514÷5=", "393÷7="),
    @("401÷9=", "977÷2="),
    @("320÷7=", "962÷7="),
    @("964÷3=", "203÷5="),
    @("867÷8=", "994÷6="),
    @("518÷2=", "218÷2="),
    @("334÷4=", "890÷2="),
    @("931÷5=", "584÷4="),
    @("618÷2=", "275÷9="),
    @("589÷6=", "728÷4="),
    @("356÷6=", "280÷5="),
    @("810÷8=", "940÷8="),
    @("943÷2=", "722÷4="),
    @("992÷4=", "820÷2="),
    @("520÷9=", "262÷7="),
    @("164÷3=", "626÷8="),
    @("112÷9=", "741÷4="),
    @("115÷8=", "681÷7="),
    @("719÷2=", "754÷5="),
    @("550÷7=", "439÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
